# Fruta / hortaliza, semanal
#
# This weekly refresh re-shuffles the price-observation columns
# (Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio,
# Unidad de comercializacion, Origen, Precio $/Kg, Kg/unidad) across the
# existing data rows, while the market/product identification columns
# (Mercado ID, Mercado, Region, Codreg, Tipo, Producto ID, Producto,
# Categoria ID, Categoria) stay put on their row.
#
# Concretely this is a permutation of rows 2..33 restricted to columns
# D, K, L, M, N, O, P, Q, R, S, T: the new value in row R for those
# columns is the old value that used to live in row Map[R].

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (1-based) that get shuffled: D=4 K=11 L=12 M=13 N=14 O=15 P=16 Q=17 R=18 S=19 T=20
$cols = @(4, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)

$firstRow = 2
$lastRow = 33

# Map: new row -> source row (the row whose old values move into the new row)
$rowMap = @{
    2  = 19
    3  = 18
    4  = 7
    5  = 16
    6  = 15
    7  = 5
    8  = 14
    9  = 17
    10 = 11
    11 = 26
    12 = 30
    13 = 29
    14 = 6
    15 = 23
    16 = 24
    17 = 21
    18 = 25
    19 = 10
    20 = 3
    21 = 4
    22 = 32
    23 = 33
    24 = 20
    25 = 22
    26 = 13
    27 = 8
    28 = 27
    29 = 12
    30 = 9
    31 = 2
    32 = 31
    33 = 28
}

# Snapshot the current (pre-shuffle) values for the columns involved, for every row,
# before writing anything back, since the permutation reads from many different rows.
# Note: .Value must be invoked as a method (parens) to actually read the cell value
# in this COM shim; a bare property access does not call the getter.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# Now write back according to the permutation.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $rowMap[$r]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $srcVals[$c]
    }
}
